$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the brand-new row 28 the same formatting as row 27 (gray fill, border,
# wrap text, Calibri) before filling it with data, so it matches the rest of
# the table instead of picking up the default style.
$ws.Range("A27:C27").Copy()
$ws.Range("A28:C28").PasteSpecial(-4122)

# Shift the existing E_QMPA / E_QMPINHABA / E_TEUR / E_PRZNTPKT block
# (rows 24-27) down by one row, working from the bottom up so no value is
# lost before it is copied.
$ws.Range("A28").Value = $ws.Range("A27").Text
$ws.Range("B28").Value = $ws.Range("B27").Text
$ws.Range("C28").Value = $ws.Range("C27").Text

$ws.Range("A27").Value = $ws.Range("A26").Text
$ws.Range("B27").Value = $ws.Range("B26").Text
$ws.Range("C27").Value = $ws.Range("C26").Text

$ws.Range("A26").Value = $ws.Range("A25").Text
$ws.Range("B26").Value = $ws.Range("B25").Text
$ws.Range("C26").Value = $ws.Range("C25").Text

$ws.Range("A25").Value = $ws.Range("A24").Text
$ws.Range("B25").Value = $ws.Range("B24").Text
$ws.Range("C25").Value = $ws.Range("C24").Text

# Row 24 now becomes E_PRZNTPKT, moved up from its old spot (row 27).
$ws.Range("A24").Value = "E_PRZNTPKT"
$ws.Range("B24").Value = "Prozentpunkte"
$ws.Range("C24").Value = "Percentage points"

# Row 28 is the brand new unit entry appended to the dictionary.
$ws.Range("A28").Value = "E_EWPKM2"
$ws.Range("B28").Value = "Einwohner/ -innen pro m² Siedlungs- und Verkehrsfläche"
$ws.Range("C28").Value = "Inhabitants per m² settlement and transport area"
